$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price values in column D (Price) to reflect latest data
$updates = @{
    "D2" = "267.88"
    "D3" = "22.93"
    "D4" = "6.319"
    "D5" = "0.06183"
    "D6" = "3.599"
    "D7" = "6.690"
    "D8" = "1.388"
    "D9" = "0.8302"
    "D10" = "0.01365"
    "D11" = "0.1600"
    "D12" = "0.08276"
    "D13" = "0.03412"
    "D14" = "0.03167"
    "D15" = "0.09281"
    "D16" = "3.897"
    "D17" = "0.001732"
    "D18" = "0.04855"
    "D19" = "0.006315"
    "D20" = "0.005373"
    "D23" = "3.770"
    "D24" = "2.327"
    "D25" = "0.3349"
    "D27" = "0.0002684"
    "D40" = "0.04663"
    "D41" = "0.006889"
    "D42" = "0.1153"
    "D43" = "0.003462"
    "D44" = "0.01214"
    "D45" = "0.00006218"
    "D47" = "0.7005"
    "D48" = "0.1753"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}

Write-Output "Done updating $($updates.Count) price cells"
